# Update cryptocurrency price/volume table with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that keeps the sheets normal (unstyled) cell format,
# used to strip any text-forcing style Excel applies automatically below.
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "60.471.96"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.601.79"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D5").Value = "'572.02"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "'142.64"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "2.625.69"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'6.49"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  -4.52%  "
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "3.067.21"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "60.515.65"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "'23.29"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "2.621.51"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'11.36"
$ws.Range("E19").Value = "  +9.01%  "
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "'346.21"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").Value = "'6.99"
$ws.Range("E22").Value = "  +7.89%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "'0.533"
$ws.Range("E24").Value = "  +13.87%  "
$ws.Range("D25").Value = "'63.27"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "'7.74"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").Value = "0.0₃0788"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "'1.84"
$ws.Range("E30").Value = "  +9.72%  "
$ws.Range("D31").Value = "'6.40"
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("D32").Value = "'0.997"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "'161.34"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "'19.49"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +4.03%  "
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  +10.17%  "
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").Value = "'1.60"
$ws.Range("E38").Value = "  +7.72%  "
$ws.Range("D39").Value = "'37.85"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").Value = "'295.69"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "'137.57"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'0.996"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "'19.74"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.94"
$ws.Range("E49").Value = "  +9.20%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0240"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'19.84"
$ws.Range("E51").Value = "  +5.46%  "

# Numeric-looking price strings must stay as text (to match the source format,
# e.g. "6.49", "0.0240"), so re-apply the normal style to drop the quote-prefix
# formatting Excel adds when forcing text entry.
$ws.Range("D5").Style = $normalStyle
$ws.Range("D6").Style = $normalStyle
$ws.Range("D7").Style = $normalStyle
$ws.Range("D10").Style = $normalStyle
$ws.Range("D16").Style = $normalStyle
$ws.Range("D19").Style = $normalStyle
$ws.Range("D21").Style = $normalStyle
$ws.Range("D22").Style = $normalStyle
$ws.Range("D24").Style = $normalStyle
$ws.Range("D25").Style = $normalStyle
$ws.Range("D26").Style = $normalStyle
$ws.Range("D28").Style = $normalStyle
$ws.Range("D30").Style = $normalStyle
$ws.Range("D31").Style = $normalStyle
$ws.Range("D32").Style = $normalStyle
$ws.Range("D33").Style = $normalStyle
$ws.Range("D34").Style = $normalStyle
$ws.Range("D36").Style = $normalStyle
$ws.Range("D38").Style = $normalStyle
$ws.Range("D39").Style = $normalStyle
$ws.Range("D42").Style = $normalStyle
$ws.Range("D43").Style = $normalStyle
$ws.Range("D44").Style = $normalStyle
$ws.Range("D47").Style = $normalStyle
$ws.Range("D49").Style = $normalStyle
$ws.Range("D50").Style = $normalStyle
$ws.Range("D51").Style = $normalStyle
